$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "한국 대학 vs. 해외 대학 (1)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/korean-uni-vs-siai-1/#utm_source=rss&utm_medium=rss&utm_campaign=korean-uni-vs-siai-1"

$ws.Range("D28").Value = "WSL2 Customize"
$ws.Range("E28").Value = "https://ropiens.tistory.com/158"
